$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price(D) and Volume(E) columns for the rows we touch,
# so Excel keeps values like "587.20" / "0.120" / "  +0.82%  " as text instead
# of converting them to numbers and dropping formatting / padding.

$textCells = @("D2","D3","D4","D5","D6","D8","D11","D12","D13","D14","D15","D17","D18","D19","D20","D21","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D36","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","E2","E3","E4","E5","E6","E7","E8","E9","E11","E12","E13","E14","E15","E16","E17","E18","E19","E20","E21","E22","E23","E24","E25","E26","E27","E28","E29","E30","E31","E32","E33","E34","E35","E36","E37","E38","E39","E40","E41","E42","E43","E44","E45","E46","E47","E48","E49","E50","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.752.04"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").Value = "2.483.64"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "587.20"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "174.05"
$ws.Range("E6").Value = "  +0.67%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  -0.46%  "

$ws.Range("E9").Value = "  +3.68%  "

$ws.Range("D11").Value = "4.95"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("D13").Value = "2.935.62"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").Value = "25.23"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").Value = "67.685.62"
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "2.456.88"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").Value = "10.78"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").Value = "7.38"
$ws.Range("E19").Value = "  -2.67%  "

$ws.Range("D20").Value = "345.80"
$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("D21").Value = "4.05"
$ws.Range("E21").Value = "  +1.31%  "

$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "70.73"
$ws.Range("E23").Value = "  +2.46%  "

$ws.Range("D24").Value = "4.17"
$ws.Range("E24").Value = "  -1.84%  "

$ws.Range("D25").Value = "1.68"
$ws.Range("E25").Value = "  -7.78%  "

$ws.Range("D26").Value = "8.82"
$ws.Range("E26").Value = "  -3.91%  "

$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("D29").Value = "0.0₃0886"
$ws.Range("E29").Value = "  -2.93%  "

$ws.Range("D30").Value = "497.81"
$ws.Range("E30").Value = "  -1.98%  "

$ws.Range("D31").Value = "7.70"
$ws.Range("E31").Value = "  -0.12%  "

$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("D33").Value = "1.76"
$ws.Range("E33").Value = "  -0.91%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("E35").Value = "  +1.26%  "

$ws.Range("D36").Value = "0.120"
$ws.Range("E36").Value = "  +1.37%  "

$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("D38").Value = "18.24"
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "1.30"
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").Value = "1.71"
$ws.Range("E41").Value = "  +1.27%  "

$ws.Range("D42").Value = "0.323"
$ws.Range("E42").Value = "  -1.91%  "

$ws.Range("D43").Value = "4.75"
$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").Value = "2.37"
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("D45").Value = "147.33"
$ws.Range("E45").Value = "  +2.70%  "

$ws.Range("D46").Value = "3.51"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").Value = "0.509"
$ws.Range("E47").Value = "  -1.41%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0252"
$ws.Range("E48").Value = "  -4.25%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.0735"
$ws.Range("E49").Value = "  -0.38%  "

$ws.Range("D50").Value = "1.55"
$ws.Range("E50").Value = "  -1.99%  "

$ws.Range("E51").Value = "  -1.66%  "
